# Apply the content edits described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("in")

# --- Text replacements in column B (Iniciativa) ---
$ws.Range("B12").Value = "Fortalecimento da Política de Patrimônio do Iepha-MG "
$ws.Range("B16").Value = "Emissão de carteiras de identificação"
$ws.Range("B52").Value = "Implantação de soluções digitais para fortalecimento da Governança Ambiental no SISEMA "
$ws.Range("B53").Value = "Diagnósticos das Cadeias Produtivas Agropecuárias"

# --- AutoFilter on the header row (adds the _FilterDatabase defined name) ---
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=in!`$A`$1:`$D`$1")
$fdb.Visible = $false

# --- Selection / view state matching the saved file ---
$ws.Range("A1:XFD1").Select()
$ws.Activate()
